$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -9
    20 = 3
    24 = -3
    25 = -1
    29 = 0
    34 = 0
    37 = 3
    38 = 1
    40 = -1
    45 = 0
    58 = 2
    59 = 5
    64 = 0
    69 = 0
    71 = 4
    75 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
